$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 1000
$ws.Range("I20").Value = 1000
$ws.Range("K20").Value = 1000
$ws.Range("M20").Value = -770

$ws.Range("H35").Value = 1000
$ws.Range("I35").Value = 1000
$ws.Range("K35").Value = 1000
$ws.Range("M35").Value = -621

$ws.Range("H62").Value = 3875.55
$ws.Range("I62").Value = 4075.4666
$ws.Range("K62").Value = 4075.4666
$ws.Range("M62").Value = -3451.4666

$ws.Range("H65").Value = 3875.55
$ws.Range("I65").Value = 4075.4666
$ws.Range("K65").Value = 20377.333
$ws.Range("M65").Value = -17257.333

$ws.Range("H86").Value = 3273657.5
$ws.Range("I86").Value = 6324.857
$ws.Range("J86").Value = 14709322
$ws.Range("K86").Value = 6324.857
$ws.Range("L86").Value = 14709322
$ws.Range("M86").Value = -5201.857
$ws.Range("N86").Value = -14711568

$ws.Range("H87").Value = 73333
$ws.Range("J87").Value = 73333
$ws.Range("L87").Value = 73333
$ws.Range("N87").Value = -75829

$ws.Range("H89").Value = 3273657.5
$ws.Range("I89").Value = 6324.857
$ws.Range("J89").Value = 14709322
$ws.Range("K89").Value = 31624.285
$ws.Range("L89").Value = 73546610
$ws.Range("M89").Value = -26008.285
$ws.Range("N89").Value = -73557842

$ws.Range("H90").Value = 73333
$ws.Range("J90").Value = 73333
$ws.Range("L90").Value = 219999
$ws.Range("N90").Value = -232479

$ws.Range("H137").Value = 1172.9706
$ws.Range("I137").Value = 1187.44
$ws.Range("K137").Value = 3562.32
$ws.Range("M137").Value = -1012.32

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2603.3647
$ws.Range("I32").Value = 2060.757
$ws.Range("K32").Value = 2060.757
$ws.Range("M32").Value = -1773.757

$ws.Range("H33").Value = 33675
$ws.Range("I33").Value = 33675
$ws.Range("K33").Value = 33675
$ws.Range("M33").Value = -33346

$ws.Range("H97").Value = 426.6875
$ws.Range("I97").Value = 459.8
$ws.Range("J97").Value = 371.5
$ws.Range("K97").Value = 459.8
$ws.Range("L97").Value = 371.5
$ws.Range("M97").Value = 36.19999999999999
$ws.Range("N97").Value = -1363.5

$ws.Range("H122").Value = 14827.714
$ws.Range("I122").Value = 5432.1113
$ws.Range("K122").Value = 16296.3339
$ws.Range("M122").Value = -13846.3339

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 879.3158
$ws.Range("J80").Value = 1008.0909
$ws.Range("L80").Value = 1008.0909
$ws.Range("N80").Value = -3004.0909

$ws.Range("H83").Value = 879.3158
$ws.Range("J83").Value = 1008.0909
$ws.Range("L83").Value = 5040.4545
$ws.Range("N83").Value = -15024.4545

$ws.Range("H86").Value = 2239.2307
$ws.Range("I86").Value = 2303.25
$ws.Range("J86").Value = 2136.8
$ws.Range("K86").Value = 2303.25
$ws.Range("L86").Value = 2136.8
$ws.Range("M86").Value = -1180.25
$ws.Range("N86").Value = -4382.8

$ws.Range("H89").Value = 2239.2307
$ws.Range("I89").Value = 2303.25
$ws.Range("J89").Value = 2136.8
$ws.Range("K89").Value = 11516.25
$ws.Range("L89").Value = 10684
$ws.Range("M89").Value = -5900.25
$ws.Range("N89").Value = -21916

$ws.Range("H94").Value = 2170.1
$ws.Range("I94").Value = 1132.5
$ws.Range("J94").Value = 6320.5
$ws.Range("K94").Value = 1132.5
$ws.Range("L94").Value = 6320.5
$ws.Range("M94").Value = -681.5
$ws.Range("N94").Value = -7222.5

$ws.Range("H96").Value = 20992.5
$ws.Range("I96").Value = 20992.5
$ws.Range("K96").Value = 20992.5
$ws.Range("M96").Value = -18246.5

$ws.Range("H107").Value = 4950.6943
$ws.Range("I107").Value = 4405.963
$ws.Range("K107").Value = 4405.963
$ws.Range("M107").Value = -2485.963

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 4002.8
$ws.Range("I10").Value = 5038
$ws.Range("J10").Value = 2450
$ws.Range("K10").Value = 5038
$ws.Range("L10").Value = 2450
$ws.Range("M10").Value = -4899
$ws.Range("N10").Value = -2728

$ws.Range("H16").Value = 166667980
$ws.Range("I16").Value = 250001000
$ws.Range("K16").Value = 250001000
$ws.Range("M16").Value = -250000713

$ws.Range("H33").Value = 4056
$ws.Range("I33").Value = 2820
$ws.Range("J33").Value = 9000
$ws.Range("K33").Value = 2820
$ws.Range("L33").Value = 9000
$ws.Range("M33").Value = -2441
$ws.Range("N33").Value = -9758

$ws.Range("H93").Value = 10500
$ws.Range("I93").Value = 10500
$ws.Range("K93").Value = 10500
$ws.Range("M93").Value = -8628

$ws.Range("H99").Value = 3064.5
$ws.Range("I99").Value = 2274.3333
$ws.Range("J99").Value = 4249.75
$ws.Range("K99").Value = 2274.3333
$ws.Range("L99").Value = 4249.75
$ws.Range("M99").Value = -776.3332999999998
$ws.Range("N99").Value = -7245.75

$ws.Range("H113").Value = 166667980
$ws.Range("I113").Value = 250001000
$ws.Range("K113").Value = 250001000
$ws.Range("M113").Value = -249998830

$ws.Range("H126").Value = 3064.5
$ws.Range("I126").Value = 2274.3333
$ws.Range("J126").Value = 4249.75
$ws.Range("K126").Value = 6822.999899999999
$ws.Range("L126").Value = 12749.25
$ws.Range("M126").Value = -4352.999899999999
$ws.Range("N126").Value = -17689.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1783.2709
$ws.Range("I131").Value = 2200.25
$ws.Range("K131").Value = 6600.75
$ws.Range("M131").Value = -1560.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 589.65
$ws.Range("I97").Value = 548.17645
$ws.Range("K97").Value = 548.17645
$ws.Range("M97").Value = -52.17645000000005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3541.0356
$ws.Range("I7").Value = 2847.647
$ws.Range("K7").Value = 2847.647
$ws.Range("M7").Value = -2735.647

$ws.Range("H16").Value = 1306.5454
$ws.Range("I16").Value = 1337.826
$ws.Range("K16").Value = 1337.826
$ws.Range("M16").Value = -1167.826

$ws.Range("H22").Value = 104299.9
$ws.Range("I22").Value = 251500
$ws.Range("J22").Value = 6166.5
$ws.Range("K22").Value = 251500
$ws.Range("L22").Value = 6166.5
$ws.Range("M22").Value = -251205
$ws.Range("N22").Value = -6756.5

$ws.Range("H25").Value = 49950
$ws.Range("I25").Value = 49950
$ws.Range("K25").Value = 49950
$ws.Range("M25").Value = -49720

$ws.Range("H27").Value = 104299.9
$ws.Range("I27").Value = 251500
$ws.Range("J27").Value = 6166.5
$ws.Range("K27").Value = 251500
$ws.Range("L27").Value = 6166.5
$ws.Range("M27").Value = -251393
$ws.Range("N27").Value = -6380.5

$ws.Range("H40").Value = 6805.857
$ws.Range("I40").Value = 7229.5557
$ws.Range("K40").Value = 7229.5557
$ws.Range("M40").Value = -7093.5557

$ws.Range("H61").Value = 1034.4
$ws.Range("I61").Value = 1034.4
$ws.Range("K61").Value = 1034.4
$ws.Range("M61").Value = -832.4000000000001

$ws.Range("H100").Value = 3937
$ws.Range("I100").Value = 2141.3333
$ws.Range("K100").Value = 2141.3333
$ws.Range("M100").Value = -1600.3333

$ws.Range("H113").Value = 1034.4
$ws.Range("I113").Value = 1034.4
$ws.Range("K113").Value = 1034.4
$ws.Range("M113").Value = 1135.6

$ws.Range("H122").Value = 5763.3335
$ws.Range("I122").Value = 5766.9287
$ws.Range("K122").Value = 17300.7861
$ws.Range("M122").Value = -14850.7861

$ws.Range("H126").Value = 3541.0356
$ws.Range("I126").Value = 2847.647
$ws.Range("K126").Value = 8542.940999999999
$ws.Range("M126").Value = -6072.940999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3705043.5
$ws.Range("I96").Value = 18519218
$ws.Range("K96").Value = 18519218
$ws.Range("M96").Value = -18517845

$ws.Range("H122").Value = 13160596
$ws.Range("I122").Value = 19233038
$ws.Range("J122").Value = 3637.5833
$ws.Range("K122").Value = 57699114
$ws.Range("L122").Value = 10912.7499
$ws.Range("M122").Value = -57696664
$ws.Range("N122").Value = -15812.7499
